# Naming Convention Guide update
#
# Applies four textual edits described by the commit diff:
#  1. Merge "Updates" + ":" into a single run "Updates:" (no text change,
#     just a run-structure cleanup).
#  2. Split "2000-2999 = Drawings" into three runs, inserting "Part " so the
#     text reads "2000-2999 = Part Drawings".
#  3. Append a new run " (only one assembly drawing per project)" right
#     after the "2000 = High level assembly drawing" run, before the
#     existing line break run.
#  4. Fix "AA,AB" -> "AA, AB" (adds a missing space) and drop the
#     proofErr gramStart/gramEnd markers that bracketed it, while keeping
#     the final result split the way the diff shows (three runs for that
#     sentence fragment).
#
# Helper: forces Word to keep a hard run boundary at a given point by
# dropping a temporary bookmark there and immediately deleting it. Word
# never merges across a point that has had a bookmark edit, but the
# bookmark itself leaves no trace once removed.
function Split-RunAt($doc, $pos) {
    $r = $doc.Range($pos, $pos)
    $doc.Bookmarks.Add("tmp_split_marker", $r) | Out-Null
    $doc.Bookmarks("tmp_split_marker").Delete()
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Updates" + ":" -> "Updates:"  (collapses the two runs into one)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Updates:", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Updates:", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. "2000-2999 = Drawings" -> "2000-2999 = Part Drawings"
#    (ending up as three runs: "2000-2999 = " / "Part " / "Drawings")
# ---------------------------------------------------------------------
$d.Content.Find.Execute("2000-2999 = Drawings", $false, $false, $false, `
    $false, $false, $true, 1, $false, "2000-2999 = Part Drawings", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("2000-2999 = Part Drawings") | Out-Null
$matchStart = $rng.Start
Split-RunAt $d ($matchStart + 12)   # after "2000-2999 = "
Split-RunAt $d ($matchStart + 17)   # after "2000-2999 = Part "

# ---------------------------------------------------------------------
# 3. "2000 = High level assembly drawing" gains a trailing run:
#    " (only one assembly drawing per project)"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("2000 = High level assembly drawing") | Out-Null
$insertPos = $rng.End
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter(" (only one assembly drawing per project)")
Split-RunAt $d $insertPos

# ---------------------------------------------------------------------
# 4. "...move to AA,AB, ..." -> "...move to AA, AB, ..." and remove the
#    proofErr gramStart/gramEnd wrapper (achieved implicitly: replacing
#    across that whole span deletes the proofErr markers because they no
#    longer sit between two surviving runs).
# ---------------------------------------------------------------------
$old4 = "24Z (If more than 26 versions are needed, move to AA,AB, "
$new4 = "24Z (If more than 26 versions are needed, move to AA, AB, "
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, `
    $true, 1, $false, $new4, 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("24Z (If more than 26 versions are needed, move to AA, AB, ") | Out-Null
$matchStart = $rng.Start
Split-RunAt $d ($matchStart + 53)   # after "...move to AA,"
Split-RunAt $d ($matchStart + 54)   # after the newly-added space
